$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1:E1").EntireColumn.Insert()
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 172000
$ws.Range("E8").Value = 227800
try {
    $ws.Range("D1:E1").EntireColumn.AutoFit()
    Write-Host "AutoFit worked"
} catch {
    Write-Host "AutoFit error: $_"
}
